$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 205, shifting existing rows 205-208 down to 206-209.
$ws.Rows("205:205").Insert()

# Populate the newly inserted row 205 with the new weekly entry.
$ws.Cells.Item(205, 1).Value = 5
$ws.Cells.Item(205, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(205, 3).Value = "Maule"
$ws.Cells.Item(205, 4).Value = 44509
$ws.Cells.Item(205, 5).Value = 7
$ws.Cells.Item(205, 6).Value = "Fruta"
$ws.Cells.Item(205, 7).Value = 100101
$ws.Cells.Item(205, 8).Value = "Berries"
$ws.Cells.Item(205, 9).Value = 100101007
$ws.Cells.Item(205, 10).Value = "Kiwi"
$ws.Cells.Item(205, 11).Value = "Hayward"
$ws.Cells.Item(205, 12).Value = "Primera"
$ws.Cells.Item(205, 13).Value = 80
$ws.Cells.Item(205, 14).Value = 17000
$ws.Cells.Item(205, 15).Value = 17000
$ws.Cells.Item(205, 16).Value = 17000
$ws.Cells.Item(205, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(205, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(205, 19).Value = 944
$ws.Cells.Item(205, 20).Value = 18

# Apply the same date-number format used by the other date cells in column D.
$ws.Cells.Item(205, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
